$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header) - subject counts for CON/STR updated
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) - updated meanEMG legmaxROM values
$ws.Range("B2").Value = -0.51334858624352486
$ws.Range("C2").Value = 5.2405407220496185
$ws.Range("D2").Value = 7.1037637969911192
$ws.Range("E2").Value = 11.218901824499449

# Row 3 (STR) - updated meanEMG legmaxROM values
$ws.Range("B3").Value = -6.5387522358148971
$ws.Range("C3").Value = 4.6333683627870021
$ws.Range("D3").Value = 14.71577472434616
$ws.Range("E3").Value = 7.02036264270017

# Reflect the narrower selection left behind by the edit
$ws.Range("B1:E3").Select()
